$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValue = @"
<rpc-reply message-id="urn:uuid:79603be3-120d-469c-b20c-6e3d28982683">
  <data>
    <network-instances>
      <network-instance>
        <name>Prueba_LxVPN</name>
        <config>
          <name>Prueba_LxVPN</name>
          <type>oc-ni-types:L3VRF</type>
        </config>
        <interfaces>
          <interface>
            <id>GigabitEthernet0/3/2</id>
            <config>
              <id>GigabitEthernet0/3/2</id>
              <interface>GigabitEthernet0/3/2</interface>
              <subinterface>0</subinterface>
            </config>
          </interface>
        </interfaces>
        <protocols>
          <protocol>
            <identifier>oc-pol-types:OSPF</identifier>
            <name>22</name>
            <config>
              <identifier>oc-pol-types:OSPF</identifier>
              <name>22</name>
            </config>
            <ospfv2>
              <global>
                <config>
                  <router-id>172.16.1.3</router-id>
                </config>
              </global>
            </ospfv2>
          </protocol>
          <protocol>
            <identifier>oc-pol-types:STATIC</identifier>
            <name>default</name>
            <config>
              <identifier>oc-pol-types:STATIC</identifier>
              <name>default</name>
            </config>
          </protocol>
          <protocol>
            <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>
            <name>default</name>
            <config>
              <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>
              <name>default</name>
            </config>
          </protocol>
        </protocols>
      </network-instance>
    </network-instances>
  </data>
</rpc-reply>
"@

$newValue = $newValue -replace "`r`n", "`n"
$newValue = $newValue + "`n"

$ws.Range("F2").Value = $newValue
